# Applies the T3 raw-lab-data cleanup edit:
#  - On "Analyseresultaten" sheet, row 6 (sample-id header row) simplify
#    the "Label(Label-suffix)" duplicated text down to just "Label-suffix"
#    for every sample column C6:X6.
#  - Fix two Analysis-name typos/labels in column A: "nitriet" -> "nitriet-N"
#    and "nitraat" -> "nitraa-N".
#  - Move the saved scroll position / active selection on that sheet
#    (cosmetic, matches the workbook's last view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analyseresultaten")

# --- Row 6: strip the "X(X-n-n)" duplication down to "X-n-n" ---
$row6 = @{
    "C6" = "CW1_EFF-1-4"
    "D6" = "INF-1-4"
    "E6" = "CW1MF01-1-5"
    "F6" = "CW1MF02-1-5"
    "G6" = "CW1MF05-1-4"
    "H6" = "CW1MF06-1-4"
    "I6" = "CW1MF10-1-5"
    "J6" = "CW1MF09-1-5"
    "K6" = "CW2MF01-1-4"
    "L6" = "CW2MF02-1-4"
    "M6" = "CW2MF05-1-4"
    "N6" = "CW2MF06-1-4"
    "O6" = "CW2MF09-1-4"
    "P6" = "CW2MF10-1-4"
    "Q6" = "CW2_EFF-1-5"
    "R6" = "CW3MF01-1-4"
    "S6" = "CW3MF02-1-4"
    "T6" = "CW3MF05-1-4"
    "U6" = "CW3MF06-1-4"
    "V6" = "CW3MF09-1-4"
    "W6" = "CW3MF10-1-4"
    "X6" = "CW3_EFF-1-4"
}

foreach ($addr in $row6.Keys) {
    $ws.Range($addr).Value = $row6[$addr]
}

# --- Column A analysis-name corrections ---
$ws.Range("A58").Value = "nitriet-N"
$ws.Range("A60").Value = "nitraa-N"

# --- Restore the view's last scroll position / selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B63").Select()

